$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.056.13"
$ws.Range("E2").Value = "  -3.78%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.746.29"
$ws.Range("E3").Value = "  -4.33%  "

# Row 4 - TetherUSD
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("E4").Value = "  -0.26%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "226.64"
$ws.Range("E5").Value = "  -3.45%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5817"
$ws.Range("E6").Value = "  -2.72%  "

# Row 7 - USDC
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.003"
$ws.Range("E7").Value = "  -0.19%  "

# Row 8 - Cardano
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2716"
$ws.Range("E8").Value = "  -0.96%  "

# Row 9 - Solana
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "23.19"
$ws.Range("E9").Value = "  -0.23%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06600"
$ws.Range("E10").Value = "  -4.82%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07515"
$ws.Range("E11").Value = "  -0.96%  "

# Row 12 - WrappedEther
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.742.48"
$ws.Range("E12").Value = "  -4.78%  "

# Row 13 - Polkadot
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.736"
$ws.Range("E13").Value = "  +0.14%  "

# Row 14 - Polygon
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6055"
$ws.Range("E14").Value = "  -2.90%  "

# Row 15 - WrappedliquidstakedEther2.0
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.982.73"
$ws.Range("E15").Value = "  -4.39%  "

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "74.30"
$ws.Range("E16").Value = "  -3.45%  "

# Row 17 - ShibaInu
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000008679"
$ws.Range("E17").Value = "  -11.24%  "

# Row 18 - WrappedBTC
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "28.048.92"
$ws.Range("E18").Value = "  -2.77%  "

# Row 19 - Uniswap
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.339"
$ws.Range("E19").Value = "  -3.63%  "

# Row 20 - Dai
$ws.Range("E20").Value = "  -0.13%  "

# Row 21 - BitcoinCash
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "205.25"
$ws.Range("E21").Value = "  -4.62%  "

# Row 22 - Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "11.27"
$ws.Range("E22").Value = "  -1.83%  "

# Row 23 - Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.641"
$ws.Range("E23").Value = "  -3.25%  "

# Row 24 - BinanceUSD
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.002"
$ws.Range("E24").Value = "  -0.24%  "

# Row 25 - Monero
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "149.92"
$ws.Range("E25").Value = "  -4.06%  "

# Row 26 - Cosmos
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.027"
$ws.Range("E26").Value = "  +1.49%  "

# Row 27 - Stellar
$ws.Range("E27").Value = "  -3.69%  "

# Row 28 - EthereumClassic
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.16"
$ws.Range("E28").Value = "  -1.63%  "

# Row 29 - Toncoin
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.390"
$ws.Range("E29").Value = "  -1.66%  "

# Row 30 - Hedera
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.06173"
$ws.Range("E30").Value = "  -4.41%  "

# Row 31 - PancakeSwap
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.391"
$ws.Range("E31").Value = "  -3.43%  "

# Row 32 - Filecoin
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.740"
$ws.Range("E32").Value = "  -2.25%  "

# Row 33 - InternetComputer(DFINITY)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.720"
$ws.Range("E33").Value = "  -1.16%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.681"
$ws.Range("E34").Value = "  -2.32%  "

# Row 35 - ARBITRUM
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.038"
$ws.Range("E35").Value = "  -4.86%  "

# Row 36 - ImmutableX
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6374"
$ws.Range("E36").Value = "  -1.15%  "

# Row 37 - HuobiToken
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.461"
$ws.Range("E37").Value = "  -3.00%  "

# Row 38 - MXToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.715"
$ws.Range("E38").Value = "  -0.96%  "

# Row 39 - VeChain
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01676"
$ws.Range("E39").Value = "  -4.33%  "

# Row 40 - FraxShare->Maker
$ws.Range("B40").Value = "Maker"
$ws.Range("C40").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.126.95"
$ws.Range("E40").Value = "  -1.11%  "

# Row 41 - Maker->FraxShare
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.220"
$ws.Range("E41").Value = "  -4.44%  "

# Row 42 - TrustWalletToken
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8735"
$ws.Range("E42").Value = "  -1.40%  "

# Row 43 - PaxDollar
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.005"
$ws.Range("E43").Value = "  +0.08%  "

# Row 44 - Quant
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "99.59"
$ws.Range("E44").Value = "  -0.31%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.896.35"
$ws.Range("E45").Value = "  -4.57%  "

# Row 46 - Aave
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "59.49"
$ws.Range("E46").Value = "  -3.26%  "

# Row 47 - RenderToken
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.585"
$ws.Range("E47").Value = "  -0.61%  "

# Row 48 - BabyDogeCoin
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000108"
$ws.Range("E48").Value = "  -4.02%  "

# Row 49 - EnergySwap
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.280"
$ws.Range("E49").Value = "  -2.15%  "

# Row 50 - Cronos
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05379"
$ws.Range("E50").Value = "  -2.10%  "

# Row 51 - Aptos
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.296"
$ws.Range("E51").Value = "  -1.29%  "

